# Update FFXIV leve-profit figures on the CRP and CUL sheets
# (values refreshed from the latest Universalis market-board snapshot).
$wb = $excel.ActiveWorkbook
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")

# Row 16 (CRP)
$wsCRP.Range("H16").Value = 1179.6
$wsCRP.Range("I16").Value = 992.73334
$wsCRP.Range("J16").Value = 1459.9
$wsCRP.Range("K16").Value = 992.73334
$wsCRP.Range("L16").Value = 1459.9
$wsCRP.Range("M16").Value = -705.73334
$wsCRP.Range("N16").Value = -2033.9

# Row 20 (CRP)
$wsCRP.Range("H20").Value = 48045.8
$wsCRP.Range("J20").Value = 48045.8
$wsCRP.Range("L20").Value = 48045.8
$wsCRP.Range("N20").Value = -48517.8

# Row 30 (CRP)
$wsCRP.Range("H30").Value = 48045.8
$wsCRP.Range("J30").Value = 48045.8
$wsCRP.Range("L30").Value = 48045.8
$wsCRP.Range("N30").Value = -48227.8

# Row 58 (CRP)
$wsCRP.Range("H58").Value = 853
$wsCRP.Range("I58").Value = 842.1142599999999
$wsCRP.Range("J58").Value = 916.5
$wsCRP.Range("K58").Value = 842.1142599999999
$wsCRP.Range("L58").Value = 916.5
$wsCRP.Range("M58").Value = -639.1142599999999
$wsCRP.Range("N58").Value = -1322.5

# Row 99 (CRP)
$wsCRP.Range("H99").Value = 2987.3333
$wsCRP.Range("I99").Value = 4133.3335
$wsCRP.Range("J99").Value = 1841.3334
$wsCRP.Range("K99").Value = 4133.3335
$wsCRP.Range("L99").Value = 1841.3334
$wsCRP.Range("M99").Value = -2635.3335
$wsCRP.Range("N99").Value = -4837.3334

# Row 107 (CRP)
$wsCRP.Range("H107").Value = 904.6286
$wsCRP.Range("I107").Value = 905.3913
$wsCRP.Range("J107").Value = 903.1667
$wsCRP.Range("K107").Value = 905.3913
$wsCRP.Range("L107").Value = 903.1667
$wsCRP.Range("M107").Value = 1014.6087
$wsCRP.Range("N107").Value = -4743.1667

# Row 108 (CRP)
$wsCRP.Range("H108").Value = 16313
$wsCRP.Range("J108").Value = 17391.25
$wsCRP.Range("L108").Value = 17391.25
$wsCRP.Range("N108").Value = -25071.25

# Row 109 (CRP)
$wsCRP.Range("H109").Value = 10999
$wsCRP.Range("J109").Value = 10999
$wsCRP.Range("L109").Value = 10999
$wsCRP.Range("N109").Value = -13079

# Row 110 (CRP)
$wsCRP.Range("H110").Value = 0
$wsCRP.Range("J110").Value = 0
$wsCRP.Range("L110").Value = 0

# Row 111 (CRP)
$wsCRP.Range("H111").Value = 34000
$wsCRP.Range("J111").Value = 34000
$wsCRP.Range("L111").Value = 34000
$wsCRP.Range("N111").Value = -42180

# Row 113 (CRP)
$wsCRP.Range("H113").Value = 1179.6
$wsCRP.Range("I113").Value = 992.73334
$wsCRP.Range("J113").Value = 1459.9
$wsCRP.Range("K113").Value = 992.73334
$wsCRP.Range("L113").Value = 1459.9
$wsCRP.Range("M113").Value = 1177.26666
$wsCRP.Range("N113").Value = -5799.9

# Row 114 (CRP)
$wsCRP.Range("H114").Value = 24900
$wsCRP.Range("J114").Value = 24900
$wsCRP.Range("L114").Value = 24900
$wsCRP.Range("N114").Value = -33578

# Row 115 (CRP)
$wsCRP.Range("H115").Value = 38000
$wsCRP.Range("J115").Value = 38000
$wsCRP.Range("L115").Value = 38000
$wsCRP.Range("N115").Value = -40350

# Row 116 (CRP)
$wsCRP.Range("H116").Value = 25000
$wsCRP.Range("J116").Value = 25000
$wsCRP.Range("L116").Value = 25000
$wsCRP.Range("N116").Value = -34178

# Row 119 (CRP)
$wsCRP.Range("H119").Value = 40699.57
$wsCRP.Range("J119").Value = 40699.57
$wsCRP.Range("L119").Value = 40699.57
$wsCRP.Range("N119").Value = -50375.57

# Row 120 (CRP)
$wsCRP.Range("H120").Value = 49982.855
$wsCRP.Range("J120").Value = 49982.855
$wsCRP.Range("L120").Value = 49982.855
$wsCRP.Range("N120").Value = -57240.855

# Row 121 (CRP)
$wsCRP.Range("H121").Value = 21396
$wsCRP.Range("J121").Value = 21396
$wsCRP.Range("L121").Value = 21396
$wsCRP.Range("N121").Value = -24016

# Row 122 (CRP)
$wsCRP.Range("H122").Value = 1758.9
$wsCRP.Range("I122").Value = 1113.6666
$wsCRP.Range("J122").Value = 2726.75
$wsCRP.Range("K122").Value = 3340.9998
$wsCRP.Range("L122").Value = 8180.25
$wsCRP.Range("M122").Value = -890.9998000000001
$wsCRP.Range("N122").Value = -13080.25

# Row 123 (CRP)
$wsCRP.Range("H123").Value = 45872.727
$wsCRP.Range("J123").Value = 45872.727
$wsCRP.Range("L123").Value = 45872.727
$wsCRP.Range("N123").Value = -55672.727

# Row 124 (CRP)
$wsCRP.Range("H124").Value = 41495
$wsCRP.Range("J124").Value = 41495
$wsCRP.Range("L124").Value = 41495
$wsCRP.Range("N124").Value = -46405

# Row 125 (CRP)
$wsCRP.Range("H125").Value = 29800
$wsCRP.Range("J125").Value = 29800
$wsCRP.Range("L125").Value = 29800
$wsCRP.Range("N125").Value = -34720

# Row 126 (CRP)
$wsCRP.Range("H126").Value = 2987.3333
$wsCRP.Range("I126").Value = 4133.3335
$wsCRP.Range("J126").Value = 1841.3334
$wsCRP.Range("K126").Value = 12400.0005
$wsCRP.Range("L126").Value = 5524.0002
$wsCRP.Range("M126").Value = -9930.000499999998
$wsCRP.Range("N126").Value = -10464.0002

# Row 127 (CRP)
$wsCRP.Range("H127").Value = 55916.668
$wsCRP.Range("J127").Value = 55916.668
$wsCRP.Range("L127").Value = 55916.668
$wsCRP.Range("N127").Value = -65836.66800000001

# Row 128 (CRP)
$wsCRP.Range("H128").Value = 48045.8
$wsCRP.Range("J128").Value = 48045.8
$wsCRP.Range("L128").Value = 48045.8
$wsCRP.Range("N128").Value = -58005.8

# Row 129 (CRP)
$wsCRP.Range("H129").Value = 49992.668
$wsCRP.Range("J129").Value = 49992.668
$wsCRP.Range("L129").Value = 49992.668
$wsCRP.Range("N129").Value = -59992.668

# Row 130 (CRP)
$wsCRP.Range("H130").Value = 43412.94
$wsCRP.Range("J130").Value = 43412.94
$wsCRP.Range("L130").Value = 43412.94
$wsCRP.Range("N130").Value = -53452.94

# Row 131 (CRP)
$wsCRP.Range("H131").Value = 34191.2
$wsCRP.Range("J131").Value = 34191.2
$wsCRP.Range("L131").Value = 34191.2
$wsCRP.Range("N131").Value = -44271.2

# Row 132 (CRP)
$wsCRP.Range("H132").Value = 2822.353
$wsCRP.Range("I132").Value = 2427.2856
$wsCRP.Range("J132").Value = 4666
$wsCRP.Range("K132").Value = 7281.8568
$wsCRP.Range("L132").Value = 13998
$wsCRP.Range("M132").Value = -4751.8568
$wsCRP.Range("N132").Value = -19058

# Row 133 (CRP)
$wsCRP.Range("H133").Value = 47750
$wsCRP.Range("J133").Value = 47750
$wsCRP.Range("L133").Value = 47750
$wsCRP.Range("N133").Value = -52810

# Row 134 (CRP)
$wsCRP.Range("H134").Value = 5900.778
$wsCRP.Range("I134").Value = 6523.091
$wsCRP.Range("K134").Value = 19569.273
$wsCRP.Range("M134").Value = -17034.273

# Row 136 (CRP)
$wsCRP.Range("H136").Value = 853
$wsCRP.Range("I136").Value = 842.1142599999999
$wsCRP.Range("J136").Value = 916.5
$wsCRP.Range("K136").Value = 2526.34278
$wsCRP.Range("L136").Value = 2749.5
$wsCRP.Range("M136").Value = 23.65722000000005
$wsCRP.Range("N136").Value = -7849.5

# Row 137 (CRP)
$wsCRP.Range("H137").Value = 39367.8
$wsCRP.Range("J137").Value = 39367.8
$wsCRP.Range("L137").Value = 39367.8
$wsCRP.Range("N137").Value = -49567.8

# Row 138 (CRP)
$wsCRP.Range("H138").Value = 39780
$wsCRP.Range("J138").Value = 39780
$wsCRP.Range("L138").Value = 39780
$wsCRP.Range("N138").Value = -50060

# Row 139 (CRP)
$wsCRP.Range("H139").Value = 41663.75
$wsCRP.Range("J139").Value = 41663.75
$wsCRP.Range("L139").Value = 41663.75
$wsCRP.Range("N139").Value = -51943.75

# Row 140 (CRP)
$wsCRP.Range("H140").Value = 71395.8
$wsCRP.Range("J140").Value = 71395.8
$wsCRP.Range("L140").Value = 71395.8
$wsCRP.Range("N140").Value = -81755.8

# Row 122 (CUL)
$wsCUL.Range("H122").Value = 54630.883
$wsCUL.Range("I122").Value = 328
$wsCUL.Range("J122").Value = 63546.285
$wsCUL.Range("K122").Value = 2952
$wsCUL.Range("L122").Value = 571916.5650000001
$wsCUL.Range("M122").Value = -502
$wsCUL.Range("N122").Value = -576816.5650000001

# Row 123 (CUL)
$wsCUL.Range("H123").Value = 5577
$wsCUL.Range("I123").Value = 1875
$wsCUL.Range("J123").Value = 6399.6665
$wsCUL.Range("K123").Value = 5625
$wsCUL.Range("L123").Value = 19198.9995
$wsCUL.Range("M123").Value = -3175
$wsCUL.Range("N123").Value = -24098.9995

# Row 124 (CUL)
$wsCUL.Range("H124").Value = 3853.75
$wsCUL.Range("I124").Value = 1357.5
$wsCUL.Range("J124").Value = 6350
$wsCUL.Range("K124").Value = 4072.5
$wsCUL.Range("L124").Value = 19050
$wsCUL.Range("M124").Value = 837.5
$wsCUL.Range("N124").Value = -28870

# Row 125 (CUL)
$wsCUL.Range("H125").Value = 2518.5186
$wsCUL.Range("I125").Value = 1000
$wsCUL.Range("J125").Value = 2640
$wsCUL.Range("K125").Value = 3000
$wsCUL.Range("L125").Value = 7920
$wsCUL.Range("M125").Value = 1920
$wsCUL.Range("N125").Value = -17760

# N110 (CRP) no longer applicable -> remove the cell entirely
$wsCRP.Range("N110").ClearContents()
